# plantilla.xlsx -> "Leer archivos de excel usando la funcion tabulate"
#
# 1. Rename the only worksheet from "Hoja1" to "sueldos".
# 2. The stray formatted-but-empty cell I4 (style only, no value) is really
#    meant to sit right after the data table (column E) instead of out at
#    column I - move it there and clear out the old spot so the used range
#    shrinks back down to A1:E5.
# 3. Column D ("CORREO") holds a long e-mail address - best-fit/auto-fit its
#    width so the text isn't clipped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename sheet -------------------------------------------------
$ws.Name = "sueldos"

# --- 2. Relocate the orphaned styled cell from I4 to E4 ---------------
$ws.Range("I4").Copy($ws.Range("E4"))
$ws.Range("I4").Clear()

# --- 3. Best-fit column D (CORREO) so the email addresses aren't cut --
$ws.Columns.Item(4).AutoFit()
